$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-16 hold pool assignments in columns A (Team 1), B (Pool 2),
# C (Pool 3); column D (Pool 4) is untouched. The edit rotates A/B/C one
# slot to the right within each row:
#   new A = old C
#   new B = old A
#   new C = old B
# Row 16 only had a value in column C, so after the rotation it only has
# a value in column A (B and C become empty).
for ($r = 2; $r -le 16; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()

    if ($c -eq $null) { $ws.Cells.Item($r, 1).Clear() } else { $ws.Cells.Item($r, 1).Value = $c }
    if ($a -eq $null) { $ws.Cells.Item($r, 2).Clear() } else { $ws.Cells.Item($r, 2).Value = $a }
    if ($b -eq $null) { $ws.Cells.Item($r, 3).Clear() } else { $ws.Cells.Item($r, 3).Value = $b }
}

# The active selection moves from B8 to A2.
$ws.Range("A2").Select() | Out-Null
